$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of results (row 3), matching the structure of row 2.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 42624.625775462962

$ws.Range("B3").Value = 65
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Random"
